$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to edit locked cells, re-protect afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidentiality footer (shared string).
$ws.Range("A40").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.
Model holdings provided as of 2021-07-07 for illustrative purposes only and are subject to change."

# Refresh Weight (D) and Percent Change (E) columns for rows 2-37.
$ws.Range("D2").Value = 0.03030033277216757
$ws.Range("E2").Value = 0
$ws.Range("D3").Value = 0.03348854223199348
$ws.Range("E3").Value = 0
$ws.Range("D4").Value = 0.03533436739438682
$ws.Range("E4").Value = 0
$ws.Range("D5").Value = 0.06874137165568085
$ws.Range("E5").Value = 0
$ws.Range("D6").Value = 0.03021542855497136
$ws.Range("E6").Value = 0
$ws.Range("D7").Value = 0.01584267809594625
$ws.Range("E7").Value = 0
$ws.Range("D8").Value = 0.03189612062532891
$ws.Range("E8").Value = 0
$ws.Range("D9").Value = 0.03101023675559667
$ws.Range("E9").Value = 0
$ws.Range("D10").Value = 0.05064704868078669
$ws.Range("E10").Value = 0
$ws.Range("D11").Value = 0.02770046839449865
$ws.Range("E11").Value = 0
$ws.Range("D12").Value = 0.01575758686505578
$ws.Range("E12").Value = 0
$ws.Range("D13").Value = 0.01619445085485829
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0.0149839112118824
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0.006607380832068808
$ws.Range("E15").Value = 0
$ws.Range("D16").Value = 0.00716187643556387
$ws.Range("E16").Value = 0
$ws.Range("D17").Value = 0.03276012389283217
$ws.Range("E17").Value = 0
$ws.Range("D18").Value = 0.0281377064116897
$ws.Range("E18").Value = 0
$ws.Range("D19").Value = 0.03068258876324477
$ws.Range("E19").Value = 0
$ws.Range("D20").Value = 0.03298734553136388
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 0.04853790823686855
$ws.Range("E21").Value = 0
$ws.Range("D22").Value = 0.02766755398430805
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 0.02985299601548623
$ws.Range("E23").Value = 0
$ws.Range("D24").Value = 0.02743827519513944
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0.01241322097051879
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 0.0131171405157314
$ws.Range("E26").Value = 0
$ws.Range("D27").Value = 0.02854913653907221
$ws.Range("E27").Value = 0
$ws.Range("D28").Value = 0.02832733829767418
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0.03115573340973466
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0.03293703884760665
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0.03096722360591577
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0.02788935222570608
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 0.0291862921954323
$ws.Range("E33").Value = 0
$ws.Range("D34").Value = 0.03068969528362684
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 0.02868977083715932
$ws.Range("E35").Value = 0
$ws.Range("D36").Value = 0.03213175788010252
$ws.Range("E36").Value = 0
$ws.Range("D37").Value = 0.9999999999999999
$ws.Range("E37").Value = 0

# Restore sheet protection (password hash is re-derived by the host; original
# protection attributes: sheet, objects, scenarios locked, format rows/cols allowed).
$ws.Protect("D382", $true, $true, $true)

